$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = "x"
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = "Hash Table"
$ws.Range("H5").Value = "String"
$ws.Range("I5").Value = "Sliding Window"

# Row 6
$ws.Range("A6").Value = "x"
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = "Hash Table"
$ws.Range("H6").Value = "String"
$ws.Range("I6").Value = "Sliding Window"

# Row 7
$ws.Range("A7").Value = "x"
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = "Hash Table"
$ws.Range("H7").Value = "String"
$ws.Range("I7").Value = "Sliding Window"

# Preserve original row heights (setting values can trigger autofit)
$ws.Rows.Item(5).RowHeight = 15.95
$ws.Rows.Item(6).RowHeight = 15.95
$ws.Rows.Item(7).RowHeight = 15.95
